# Apply symbol-list update (commit: 'Updated symbol list on Tue Feb 14 07:40:00 UTC 2023 with GitHub Actions')
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (Coin name / Link columns) -------------------------
$textUpdates = @(
    @{Cell='B15'; Value='TigerCash'},
    @{Cell='C15'; Value='https://coinranking.com/coin/6hIn06L2+tigercash-tch'},
    @{Cell='B16'; Value='LEO'},
    @{Cell='C16'; Value='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'},
    @{Cell='B17'; Value='GateToken'},
    @{Cell='C17'; Value='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'},
    @{Cell='B18'; Value='BitpandaEcosystemToken'},
    @{Cell='C18'; Value='https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'},
    @{Cell='B19'; Value='MCDex'},
    @{Cell='C19'; Value='https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'},
    @{Cell='B20'; Value='ProBitToken'},
    @{Cell='C20'; Value='https://coinranking.com/coin/lQP4d6T2+probittoken-prob'},
    @{Cell='B21'; Value='ZBToken'},
    @{Cell='C21'; Value='https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'},
    @{Cell='B22'; Value='CoinExToken'},
    @{Cell='C22'; Value='https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'}
)
foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- Numeric-looking cells (Price / Volume columns) -----------------------
# Forced to Text format before the write so Excel keeps the exact literal
# string (matches the source inlineStr cells) instead of coercing to a
# number/percentage; format is cleared again afterwards to avoid leaving
# the cells tagged as Text-formatted.
$numericUpdates = @(
    @{Cell='D2'; Value='293.19'},
    @{Cell='E2'; Value='-6.77%'},
    @{Cell='D3'; Value='40.61'},
    @{Cell='E3'; Value='-0.83%'},
    @{Cell='D4'; Value='5.027'},
    @{Cell='E4'; Value='-2.22%'},
    @{Cell='D5'; Value='0.07328'},
    @{Cell='E5'; Value='-3.63%'},
    @{Cell='D6'; Value='1.539'},
    @{Cell='E6'; Value='-8.50%'},
    @{Cell='D7'; Value='0.9278'},
    @{Cell='E7'; Value='-0.34%'},
    @{Cell='D9'; Value='0.1164'},
    @{Cell='E9'; Value='-3.16%'},
    @{Cell='E10'; Value='-4.68%'},
    @{Cell='D11'; Value='0.04346'},
    @{Cell='E11'; Value='5.21%'},
    @{Cell='D12'; Value='0.08723'},
    @{Cell='E12'; Value='-3.71%'},
    @{Cell='E13'; Value='0.05%'},
    @{Cell='E14'; Value='-0.35%'},
    @{Cell='D15'; Value='0.006043'},
    @{Cell='E15'; Value='3.11%'},
    @{Cell='D16'; Value='3.336'},
    @{Cell='E16'; Value='0.04%'},
    @{Cell='D17'; Value='4.284'},
    @{Cell='E17'; Value='-0.84%'},
    @{Cell='D18'; Value='0.3289'},
    @{Cell='E18'; Value='-1.84%'},
    @{Cell='D19'; Value='7.971'},
    @{Cell='E19'; Value='4.92%'},
    @{Cell='D20'; Value='0.1400'},
    @{Cell='E20'; Value='3.64%'},
    @{Cell='D21'; Value='0.2742'},
    @{Cell='E21'; Value='-3.25%'},
    @{Cell='D22'; Value='0.03942'},
    @{Cell='E22'; Value='-1.10%'},
    @{Cell='E23'; Value='-1.36%'},
    @{Cell='D24'; Value='0.003788'},
    @{Cell='E24'; Value='-6.79%'},
    @{Cell='E25'; Value='-5.08%'},
    @{Cell='D26'; Value='0.0003724'},
    @{Cell='E26'; Value='22.60%'},
    @{Cell='D38'; Value='0.02309'},
    @{Cell='E38'; Value='-5.15%'},
    @{Cell='D39'; Value='0.05076'},
    @{Cell='E39'; Value='-2.18%'},
    @{Cell='D40'; Value='0.006115'},
    @{Cell='E40'; Value='85.31%'},
    @{Cell='D41'; Value='0.007855'},
    @{Cell='E41'; Value='2.02%'},
    @{Cell='D42'; Value='0.1287'},
    @{Cell='E42'; Value='-1.25%'},
    @{Cell='D43'; Value='0.007356'},
    @{Cell='E43'; Value='-3.03%'},
    @{Cell='D44'; Value='0.007268'},
    @{Cell='E44'; Value='-14.30%'},
    @{Cell='D45'; Value='0.3195'},
    @{Cell='E45'; Value='-5.87%'},
    @{Cell='D46'; Value='0.00006287'},
    @{Cell='E46'; Value='-4.52%'},
    @{Cell='E47'; Value='0.13%'},
    @{Cell='D48'; Value='0.03379'},
    @{Cell='E48'; Value='-87.66%'},
    @{Cell='D49'; Value='0.00002101'},
    @{Cell='E49'; Value='0.13%'},
    @{Cell='D50'; Value='0.0002001'},
    @{Cell='E50'; Value='0.13%'}
)
foreach ($u in $numericUpdates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = '@'
    $rng.Value = $u.Value
}

# Clear the temporary Text number-format from the touched rows' D:E cells
# (restores 'General' formatting to match the original file).
$ws.Range('D2:E7').ClearFormats()
$ws.Range('D9:E26').ClearFormats()
$ws.Range('D38:E50').ClearFormats()
